$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet "Feuil1" to "p4" ---
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws1.Name = "p4"

# --- p4: drop the "Classe" column (D) and the second data row (test@gmail.com),
#     which moves to its own new sheet "p5" ---

# Remove the hyperlinks on p4 entirely (cannot selectively remove a single
# hyperlink via this host), then re-create only the one we keep (A2).
$ws1.Hyperlinks.Delete()

$ws1.Range("B3").ClearContents()
$ws1.Range("C3").ClearContents()
$ws1.Range("D1").ClearContents()
$ws1.Range("D2").ClearContents()
$ws1.Range("D3").ClearContents()
$ws1.Range("A3").ClearContents()

$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:louiscarlier123@gmail.com") | Out-Null
$ws1.Range("A2").Style = "Lien hypertexte"

# D1 keeps a light fill but no value/header anymore
$ws1.Range("D1").Interior.ThemeColor = 1

$ws1.Range("A15").Select()

# --- Create the new sheet "p5" with the rows that used to live in p4 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "p5"

$ws2.Range("A1").Value = "Email"
$ws2.Range("B1").Value = "Nom"
$ws2.Range("C1").Value = "Prenom"

$ws2.Range("A2").Value = "test@gmail.com"
$ws2.Range("B2").Value = "testname"
$ws2.Range("C2").Value = "testfirstname"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:test@gmail.com") | Out-Null
$ws2.Range("A2").Style = "Lien hypertexte"

$ws2.Range("D1").Interior.ThemeColor = 1

$ws2.Range("B17").Select()
